$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Toyota Yaris"
$ws.Range("C1").Value = "Mazda MX-30"
$ws.Range("D1").Value = "Honda JAZZ"
$ws.Range("E1").Value = "Land Rover Defender"
$ws.Range("F1").Value = "SEAT Leon"
$ws.Range("G1").Value = "KIA Sorento"
$ws.Range("H1").Value = "Honda e"
$ws.Range("I1").Value = "Hyundai i10"
$ws.Range("J1").Value = "ISUZU D-Max Crew Cab"
$ws.Range("K1").Value = "Audi A3"
